$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; existing rows 15-55 shift down to 16-56.
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = 44883
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100103
$ws.Cells.Item(15, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(15, 9).Value = 100103004
$ws.Cells.Item(15, 10).Value = "Durazno"
$ws.Cells.Item(15, 11).Value = "Florida King"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 300
$ws.Cells.Item(15, 14).Value = 24000
$ws.Cells.Item(15, 15).Value = 25000
$ws.Cells.Item(15, 16).Value = 24500
$ws.Cells.Item(15, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(15, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(15, 19).Value = 1361
$ws.Cells.Item(15, 20).Value = 18
